$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title shape "A" + " " + "slide" -> "A slide" (merge runs into one)
$titleShape = $s.Shapes.Item("Title 1")
$titleShape.TextFrame.TextRange.Text = "__TEMP__"
$titleShape.TextFrame.TextRange.Text = "A slide"

# Table cell "a" + " " + "table" -> "a table" (merge runs into one)
$tableShape = $s.Shapes.Item("Content Placeholder 5")
$cell = $tableShape.Table.Cell(1, 2)
$cell.Shape.TextFrame.TextRange.Text = "__TEMP__"
$cell.Shape.TextFrame.TextRange.Text = "a table"

# TextBox "Plus" + " " + "an" + " " + "image" -> "Plus an image" (merge runs into one)
$textBoxShape = $s.Shapes.Item("TextBox 3")
$textBoxShape.TextFrame.TextRange.Text = "__TEMP__"
$textBoxShape.TextFrame.TextRange.Text = "Plus an image"
